# Integrated model stats from text-based sim
# Adds a "Radius" column (L) to the Templar Models sheet, with a bold,
# left/right-bordered header and a value of 12 for every populated model row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Templar Models")
# ($wb.ActiveSheet resolves to the same sheet here, since it is the tab
# that was selected when the workbook was last saved.)

# Header cell: label + bold font + thin left/right border
$ws.Range("L2").Value = "Radius"
$ws.Range("L2").Font.Bold = $true
$ws.Range("L2").Borders.Item(7).LineStyle = 1
$ws.Range("L2").Borders.Item(10).LineStyle = 1

# Radius values for each populated model row
$ws.Range("L3").Value = 12
$ws.Range("L4").Value = 12
$ws.Range("L5").Value = 12
$ws.Range("L8").Value = 12
$ws.Range("L9").Value = 12
$ws.Range("L10").Value = 12

# Move the active selection to the new header cell
$ws.Range("L2").Select()
